# Add two new columns "I0" (col I) and "IF" (col J) to the sheet,
# matching the header style used by the existing header row, and populate
# the per-row values for rows 2-60.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers: copy the style from the existing header cell (H1) onto I1/J1.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("J1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Per-row data for the new columns.
$data = @(
    @{ Row = 2; I = 7; J = 7 },
    @{ Row = 3; I = 8; J = 8 },
    @{ Row = 4; I = 6; J = 6 },
    @{ Row = 5; I = 7; J = 7 },
    @{ Row = 6; I = 7; J = 7 },
    @{ Row = 7; I = 8; J = 8 },
    @{ Row = 8; I = 7; J = 7 },
    @{ Row = 9; I = 6; J = 6 },
    @{ Row = 10; I = 6; J = 7 },
    @{ Row = 11; I = 8; J = 8 },
    @{ Row = 12; I = 8; J = 8 },
    @{ Row = 13; I = 7; J = 7 },
    @{ Row = 14; I = 8; J = 8 },
    @{ Row = 15; I = 8; J = 8 },
    @{ Row = 16; I = 8; J = 8 },
    @{ Row = 17; I = 8; J = 8 },
    @{ Row = 18; I = 8; J = 8 },
    @{ Row = 19; I = 7; J = 7 },
    @{ Row = 20; I = 8; J = 8 },
    @{ Row = 21; I = 7; J = 7 },
    @{ Row = 22; I = 7; J = 7 },
    @{ Row = 23; I = 7; J = 8 },
    @{ Row = 24; I = 8; J = 8 },
    @{ Row = 25; I = 8; J = 8 },
    @{ Row = 26; I = 8; J = 8 },
    @{ Row = 27; I = 8; J = 8 },
    @{ Row = 28; I = 8; J = 9 },
    @{ Row = 29; I = 8; J = 8 },
    @{ Row = 30; I = 8; J = 8 },
    @{ Row = 31; I = 8; J = 8 },
    @{ Row = 32; I = 8; J = 8 },
    @{ Row = 33; I = 8; J = 8 },
    @{ Row = 34; I = 9; J = 9 },
    @{ Row = 35; I = 7; J = 8 },
    @{ Row = 36; I = 7; J = 8 },
    @{ Row = 37; I = 9; J = 9 },
    @{ Row = 38; I = 8; J = 8 },
    @{ Row = 39; I = 8; J = 8 },
    @{ Row = 40; I = 7; J = 7 },
    @{ Row = 41; I = 7; J = 8 },
    @{ Row = 42; I = 8; J = 8 },
    @{ Row = 43; I = 8; J = 8 },
    @{ Row = 44; I = 9; J = 10 },
    @{ Row = 45; I = 7; J = 7 },
    @{ Row = 46; I = 5; J = 6 },
    @{ Row = 47; I = 6; J = 6 },
    @{ Row = 48; I = 9; J = 9 },
    @{ Row = 49; I = 7; J = 7 },
    @{ Row = 50; I = 7; J = 7 },
    @{ Row = 51; I = 6; J = 6 },
    @{ Row = 52; I = 7; J = 7 },
    @{ Row = 53; I = 3; J = 3 },
    @{ Row = 54; I = 7; J = 7 },
    @{ Row = 55; I = 6; J = 7 },
    @{ Row = 56; I = 6; J = 6 },
    @{ Row = 57; I = 6; J = 6 },
    @{ Row = 58; I = 6; J = 6 },
    @{ Row = 59; I = 5; J = 5 },
    @{ Row = 60; I = 6; J = 6 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 9).Value = $item.I
    $ws.Cells.Item($r, 10).Value = $item.J
}
